# Update test-scenario names in column A to match the revised wording.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "Viewing list of posts by category"
$ws.Range("A7").Value = "Leaving a comment with the author more than 60 characters"
$ws.Range("A8").Value = "Leaving a comment without an author"
$ws.Range("A9").Value = "Leaving a comment without a body"

# Row heights re-flowed (content wrapped differently after the wording edit).
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 105
$ws.Rows.Item(9).RowHeight = 105

# Move the viewport / selection to the newly edited rows.
$null = $ws.Range("A9").Select()
